$wb = $excel.ActiveWorkbook

# ALC row 18 (diff @@ -1526)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1431.7
$ws.Range("I18").Value = 1335.2222
$ws.Range("K18").Value = 1335.2222
$ws.Range("M18").Value = -1051.2222

# ALC row 33 (diff @@ -2267)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7972.7
$ws.Range("I33").Value = 11210.714
$ws.Range("J33").Value = 417.33334
$ws.Range("K33").Value = 11210.714
$ws.Range("L33").Value = 417.33334
$ws.Range("M33").Value = -10981.714
$ws.Range("N33").Value = -875.33334

# ALC row 58 (diff @@ -3534)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 12763.5
$ws.Range("I58").Value = 12763.5
$ws.Range("K58").Value = 38290.5
$ws.Range("M58").Value = -38140.5

# ALC row 62 (diff @@ -3733)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 84446500
$ws.Range("I62").Value = 105557210
$ws.Range("K62").Value = 105557210
$ws.Range("M62").Value = -105556586

# ALC row 65 (diff @@ -3880)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 84446500
$ws.Range("I65").Value = 105557210
$ws.Range("K65").Value = 527786050
$ws.Range("M65").Value = -527782930

# ALC row 138 (diff @@ -7535)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2864.375
$ws.Range("J138").Value = 3365.054
$ws.Range("L138").Value = 10095.162
$ws.Range("N138").Value = -20375.162

# ARM row 4 (diff @@ -7932)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 542.7
$ws.Range("I4").Value = 530.625
$ws.Range("K4").Value = 530.625
$ws.Range("M4").Value = -414.625

# ARM row 32 (diff @@ -9337)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22474.865
$ws.Range("I32").Value = 22327.314
$ws.Range("K32").Value = 22327.314
$ws.Range("M32").Value = -22040.314

# ARM row 45 (diff @@ -9989)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7186.222
$ws.Range("I45").Value = 8452.714
$ws.Range("J45").Value = 2753.5
$ws.Range("K45").Value = 8452.714
$ws.Range("L45").Value = 2753.5
$ws.Range("M45").Value = -8075.714
$ws.Range("N45").Value = -3507.5

# ARM row 61 (diff @@ -10764)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4440.852
$ws.Range("I61").Value = 2243.3635
$ws.Range("K61").Value = 2243.3635
$ws.Range("M61").Value = -2031.3635

# ARM row 110 (diff @@ -13150)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 16667810
$ws.Range("J110").Value = 1666.6666
$ws.Range("L110").Value = 1666.6666
$ws.Range("N110").Value = -5756.6666

# ARM row 132 (diff @@ -14213)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3486.074
$ws.Range("I132").Value = 1982.909
$ws.Range("K132").Value = 5948.727000000001
$ws.Range("M132").Value = -3418.727000000001

# ARM row 136 (diff @@ -14412)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4440.852
$ws.Range("I136").Value = 2243.3635
$ws.Range("K136").Value = 6730.0905
$ws.Range("M136").Value = -4180.0905

# BSM row 22 (diff @@ -15798)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 693.5
$ws.Range("I22").Value = 693.5
$ws.Range("K22").Value = 693.5
$ws.Range("M22").Value = -520.5

# BSM row 80 (diff @@ -18616)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6385.4287
$ws.Range("J80").Value = 10549.75
$ws.Range("L80").Value = 10549.75
$ws.Range("N80").Value = -12545.75

# BSM row 83 (diff @@ -18769)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 6385.4287
$ws.Range("J83").Value = 10549.75
$ws.Range("L83").Value = 52748.75
$ws.Range("N83").Value = -62732.75

# BSM row 134 (diff @@ -21271)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7425.6875
$ws.Range("I134").Value = 2607.4
$ws.Range("J134").Value = 9615.817999999999
$ws.Range("K134").Value = 7822.200000000001
$ws.Range("L134").Value = 28847.454
$ws.Range("M134").Value = -5287.200000000001
$ws.Range("N134").Value = -33917.454

# CRP row 16 (diff @@ -22479)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6637.625
$ws.Range("I16").Value = 1037
$ws.Range("K16").Value = 1037
$ws.Range("M16").Value = -750

# CRP row 31 (diff @@ -23208)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41671924
$ws.Range("I31").Value = 90910490
$ws.Range("K31").Value = 90910490
$ws.Range("M31").Value = -90910195

# CRP row 34 (diff @@ -23361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 41671924
$ws.Range("I34").Value = 90910490
$ws.Range("K34").Value = 90910490
$ws.Range("M34").Value = -90910288

# CRP row 113 (diff @@ -27280)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 6637.625
$ws.Range("I113").Value = 1037
$ws.Range("K113").Value = 1037
$ws.Range("M113").Value = 1133

# CRP row 134 (diff @@ -28315)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4431.1284
$ws.Range("I134").Value = 4082.5293
$ws.Range("K134").Value = 12247.5879
$ws.Range("M134").Value = -9712.5879

# CUL row 23 (diff @@ -29872)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 83470.836
$ws.Range("I23").Value = 250055
$ws.Range("J23").Value = 178.75
$ws.Range("K23").Value = 750165
$ws.Range("L23").Value = 536.25
$ws.Range("M23").Value = -749930
$ws.Range("N23").Value = -1006.25

# CUL row 33 (diff @@ -30386)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 111.666664
$ws.Range("I33").Value = 99
$ws.Range("K33").Value = 594
$ws.Range("M33").Value = -311

# CUL row 51 (diff @@ -31298)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3918.889
$ws.Range("J51").Value = 4928.5713
$ws.Range("L51").Value = 14785.7139
$ws.Range("N51").Value = -15705.7139

# CUL row 125 (diff @@ -34978)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2998.5
$ws.Range("I125").Value = 2998.5
$ws.Range("K125").Value = 8995.5
$ws.Range("M125").Value = -4075.5

# CUL row 132 (diff @@ -35324)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 8608.166999999999
$ws.Range("I132").Value = 8608.166999999999
$ws.Range("K132").Value = 77473.503
$ws.Range("M132").Value = -74943.503

# GSM row 122 (diff @@ -41794)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5321.0938
$ws.Range("I122").Value = 5842.857
$ws.Range("J122").Value = 1668.75
$ws.Range("K122").Value = 17528.571
$ws.Range("L122").Value = 5006.25
$ws.Range("M122").Value = -15078.571
$ws.Range("N122").Value = -9906.25

# GSM row 126 (diff @@ -41987)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3447.9524
$ws.Range("I126").Value = 2759
$ws.Range("J126").Value = 4366.5557
$ws.Range("K126").Value = 8277
$ws.Range("L126").Value = 13099.6671
$ws.Range("M126").Value = -5807
$ws.Range("N126").Value = -18039.6671

# LTW row 9 (diff @@ -43199)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 588.1818
$ws.Range("I9").Value = 145.6
$ws.Range("J9").Value = 957
$ws.Range("K9").Value = 145.6
$ws.Range("L9").Value = 957
$ws.Range("M9").Value = 78.40000000000001
$ws.Range("N9").Value = -1405

# LTW row 19 (diff @@ -43686)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 3097.889
$ws.Range("J19").Value = 3280.1667
$ws.Range("L19").Value = 3280.1667
$ws.Range("N19").Value = -3620.1667

# LTW row 61 (diff @@ -45753)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1299.8334
$ws.Range("I61").Value = 934.8570999999999
$ws.Range("K61").Value = 934.8570999999999
$ws.Range("M61").Value = -732.8570999999999

# LTW row 100 (diff @@ -47655)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 15627579
$ws.Range("I100").Value = 125000950
$ws.Range("K100").Value = 125000950
$ws.Range("M100").Value = -125000409

# LTW row 113 (diff @@ -48292)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1299.8334
$ws.Range("I113").Value = 934.8570999999999
$ws.Range("K113").Value = 934.8570999999999
$ws.Range("M113").Value = 1235.1429

# LTW row 132 (diff @@ -49217)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4469.9375
$ws.Range("I132").Value = 3642.818
$ws.Range("J132").Value = 6289.6
$ws.Range("K132").Value = 10928.454
$ws.Range("L132").Value = 18868.8
$ws.Range("M132").Value = -8398.454000000002
$ws.Range("N132").Value = -23928.8

# LTW row 136 (diff @@ -49419)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5165.6772
$ws.Range("J136").Value = 8656.666999999999
$ws.Range("L136").Value = 25970.001
$ws.Range("N136").Value = -31070.001

# WVR row 54 (diff @@ -52370)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 115038.5
$ws.Range("J54").Value = 30077
$ws.Range("L54").Value = 30077
$ws.Range("N54").Value = -31117

# WVR row 123 (diff @@ -55727)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# WVR row 132 (diff @@ -56162)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6727.0386
$ws.Range("I132").Value = 4066.6667
$ws.Range("J132").Value = 9007.357
$ws.Range("K132").Value = 12200.0001
$ws.Range("L132").Value = 27022.071
$ws.Range("M132").Value = -9670.000100000001
$ws.Range("N132").Value = -32082.071
